$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching style of existing header G1
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add the corresponding data value in H2
$ws.Range("H2").Value = 1
